$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8951720595359802
$ws.Range("B1").Value = 1.769679069519043
$ws.Range("C1").Value = 4.225603103637695
$ws.Range("D1").Value = 3.523576498031616
$ws.Range("E1").Value = 1.23878002166748
